{"js": "// Update benchmark stats table: refresh summary rows (0-11) and collapse the\n// trailing per-iteration breakdown rows (43-45) down to their headline value,\n// matching the values that rows 0-2 now hold above.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Map of 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"817\",\n  5: \"0.00064\",\n  6: \"0.00020\",\n  7: \"0.00005\",\n  8: \"0.00031\",\n  9: \"0.00040\",\n  10: \"0.00047\",\n  11: \"0.16301\",\n  43: \"99.88\",\n  44: \"0.16\",\n  45: \"131\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(parseInt(rowIndex, 10), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table: refresh summary rows (1-12) and collapse the\n# trailing per-iteration breakdown rows (44-46) down to their headline value,\n# matching the values that rows 1-3 now hold above.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"817\"\n$t.Cell(6, 1).Range.Text = \"0.00064\"\n$t.Cell(7, 1).Range.Text = \"0.00020\"\n$t.Cell(8, 1).Range.Text = \"0.00005\"\n$t.Cell(9, 1).Range.Text = \"0.00031\"\n$t.Cell(10, 1).Range.Text = \"0.00040\"\n$t.Cell(11, 1).Range.Text = \"0.00047\"\n$t.Cell(12, 1).Range.Text = \"0.16301\"\n$t.Cell(44, 1).Range.Text = \"99.88\"\n$t.Cell(45, 1).Range.Text = \"0.16\"\n$t.Cell(46, 1).Range.Text = \"131\"\n"}
